$wb = $excel.ActiveWorkbook

# ===== Sheet cell values =====

# --- Overview ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A1").Value = "File Name"
$ws.Range("B1").Value = "zh-cn"
$ws.Range("C1").Value = "de-de"
$ws.Range("A2").Value = "ffffcafed8d6-ffc7-4f2e-b2a7-3fd220cb73fa.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("A3").Value = "ffffff5d935b58-0829-45f7-9e5b-8e53a6fa3f9f.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("A4").Value = "f29e047e-44c6-4de6-bfaa-fb03f56fc80b.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("A5").Value = ".localization-config"
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("C5").Value = "Not to be localized"

# --- zh-cn ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A1").Value = "Source File Name"
$ws.Range("B1").Value = "Status"
$ws.Range("C1").Value = "Latest Handoff File"
$ws.Range("D1").Value = "Latest Handoff Datetime"
$ws.Range("E1").Value = "Latest Target File"
$ws.Range("F1").Value = "Latest Handback File"
$ws.Range("G1").Value = "Latest Handback DateTime"
$ws.Range("H1").Value = "Handoff Reason"
$ws.Range("I1").Value = "Dependency From"
$ws.Range("A2").Value = "ffffcafed8d6-ffc7-4f2e-b2a7-3fd220cb73fa.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.zh-cn.xlf"
$ws.Range("D2").Value = "2016-03-10 19:06:11"
$ws.Range("E2").Value = "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.md"
$ws.Range("F2").Value = "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.zh-cn.xlf"
$ws.Range("G2").Value = "2016-03-10 19:06:40"
$ws.Range("H2").Value = "Include"
$ws.Range("A3").Value = "ffffff5d935b58-0829-45f7-9e5b-8e53a6fa3f9f.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.zh-cn.xlf"
$ws.Range("D3").Value = "2016-03-10 19:06:11"
$ws.Range("E3").Value = "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.md"
$ws.Range("F3").Value = "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.zh-cn.xlf"
$ws.Range("G3").Value = "2016-03-10 19:06:40"
$ws.Range("H3").Value = "Include"
$ws.Range("A4").Value = "f29e047e-44c6-4de6-bfaa-fb03f56fc80b.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "f29e047e-44c6-4de6-bfaa-fb03f56fc80b.543a0a31b72e48fbef1d9e1c94867f8b9fb36ec5.zh-cn.xlf"
$ws.Range("D4").Value = "2016-03-10 19:09:50"
$ws.Range("E4").Value = "f29e047e-44c6-4de6-bfaa-fb03f56fc80b.md"
$ws.Range("F4").Value = "f29e047e-44c6-4de6-bfaa-fb03f56fc80b.543a0a31b72e48fbef1d9e1c94867f8b9fb36ec5.zh-cn.xlf"
$ws.Range("G4").Value = "2016-03-10 19:09:10"
$ws.Range("H4").Value = "Include"
$ws.Range("A5").Value = ".localization-config"
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("D5").Value = "0001-01-01 00:00:00"
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Ignored"

# --- de-de ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A1").Value = "Source File Name"
$ws.Range("B1").Value = "Status"
$ws.Range("C1").Value = "Latest Handoff File"
$ws.Range("D1").Value = "Latest Handoff Datetime"
$ws.Range("E1").Value = "Latest Target File"
$ws.Range("F1").Value = "Latest Handback File"
$ws.Range("G1").Value = "Latest Handback DateTime"
$ws.Range("H1").Value = "Handoff Reason"
$ws.Range("I1").Value = "Dependency From"
$ws.Range("A2").Value = "ffffcafed8d6-ffc7-4f2e-b2a7-3fd220cb73fa.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.de-de.xlf"
$ws.Range("D2").Value = "2016-03-10 19:06:16"
$ws.Range("E2").Value = "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.md"
$ws.Range("F2").Value = "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.de-de.xlf"
$ws.Range("G2").Value = "2016-03-10 19:06:51"
$ws.Range("H2").Value = "Include"
$ws.Range("A3").Value = "ffffff5d935b58-0829-45f7-9e5b-8e53a6fa3f9f.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.de-de.xlf"
$ws.Range("D3").Value = "2016-03-10 19:06:16"
$ws.Range("E3").Value = "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.md"
$ws.Range("F3").Value = "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.de-de.xlf"
$ws.Range("G3").Value = "2016-03-10 19:06:51"
$ws.Range("H3").Value = "Include"
$ws.Range("A4").Value = "f29e047e-44c6-4de6-bfaa-fb03f56fc80b.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "f29e047e-44c6-4de6-bfaa-fb03f56fc80b.543a0a31b72e48fbef1d9e1c94867f8b9fb36ec5.de-de.xlf"
$ws.Range("D4").Value = "2016-03-10 19:09:54"
$ws.Range("E4").Value = "f29e047e-44c6-4de6-bfaa-fb03f56fc80b.md"
$ws.Range("F4").Value = "f29e047e-44c6-4de6-bfaa-fb03f56fc80b.543a0a31b72e48fbef1d9e1c94867f8b9fb36ec5.de-de.xlf"
$ws.Range("G4").Value = "2016-03-10 19:09:21"
$ws.Range("H4").Value = "Include"
$ws.Range("A5").Value = ".localization-config"
$ws.Range("B5").Value = "Not to be localized"
$ws.Range("D5").Value = "0001-01-01 00:00:00"
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Ignored"

# ===== Rebuild hyperlinks (delete + re-add in target order) =====

# --- Overview hyperlinks ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/68472668c37039b194656fc7aed7779ca3025a17/e2e/f29e047e-44c6-4de6-bfaa-fb03f56fc80b.md", [Type]::Missing, [Type]::Missing, "ffffcafed8d6-ffc7-4f2e-b2a7-3fd220cb73fa.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/68472668c37039b194656fc7aed7779ca3025a17/e2e/ffffcafed8d6-ffc7-4f2e-b2a7-3fd220cb73fa.md", [Type]::Missing, [Type]::Missing, "ffffff5d935b58-0829-45f7-9e5b-8e53a6fa3f9f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/68472668c37039b194656fc7aed7779ca3025a17/e2e/ffffff5d935b58-0829-45f7-9e5b-8e53a6fa3f9f.md", [Type]::Missing, [Type]::Missing, "f29e047e-44c6-4de6-bfaa-fb03f56fc80b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/68472668c37039b194656fc7aed7779ca3025a17/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

# --- zh-cn hyperlinks ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/68472668c37039b194656fc7aed7779ca3025a17/e2e/f29e047e-44c6-4de6-bfaa-fb03f56fc80b.md", [Type]::Missing, [Type]::Missing, "ffffcafed8d6-ffc7-4f2e-b2a7-3fd220cb73fa.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/acbd9f662ef7753b136202c2f7fa01f40e20af2b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f29e047e-44c6-4de6-bfaa-fb03f56fc80b.543a0a31b72e48fbef1d9e1c94867f8b9fb36ec5.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/a1a4c4a4d66af537b67b9ee88943ab45d8e3e926/e2e/f29e047e-44c6-4de6-bfaa-fb03f56fc80b.md", [Type]::Missing, [Type]::Missing, "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9a185fa12125651b0257d91400c949a1f5ada9c3/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f29e047e-44c6-4de6-bfaa-fb03f56fc80b.543a0a31b72e48fbef1d9e1c94867f8b9fb36ec5.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/68472668c37039b194656fc7aed7779ca3025a17/e2e/ffffcafed8d6-ffc7-4f2e-b2a7-3fd220cb73fa.md", [Type]::Missing, [Type]::Missing, "ffffff5d935b58-0829-45f7-9e5b-8e53a6fa3f9f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/969fb623a3612c61acf3c447c540bb0a8d227084/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/1cbc950f2376e7d8a3f685e557f93e0f34a72896/e2e/a28591f3-d3b4-45d5-86e2-cebe1e59fd36.md", [Type]::Missing, [Type]::Missing, "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a7e9896b6a878915ef5ab977f5c038b3499efbb1/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/68472668c37039b194656fc7aed7779ca3025a17/e2e/ffffff5d935b58-0829-45f7-9e5b-8e53a6fa3f9f.md", [Type]::Missing, [Type]::Missing, "f29e047e-44c6-4de6-bfaa-fb03f56fc80b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/969fb623a3612c61acf3c447c540bb0a8d227084/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "f29e047e-44c6-4de6-bfaa-fb03f56fc80b.543a0a31b72e48fbef1d9e1c94867f8b9fb36ec5.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/1cbc950f2376e7d8a3f685e557f93e0f34a72896/e2e/a28591f3-d3b4-45d5-86e2-cebe1e59fd36.md", [Type]::Missing, [Type]::Missing, "f29e047e-44c6-4de6-bfaa-fb03f56fc80b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a7e9896b6a878915ef5ab977f5c038b3499efbb1/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "f29e047e-44c6-4de6-bfaa-fb03f56fc80b.543a0a31b72e48fbef1d9e1c94867f8b9fb36ec5.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/68472668c37039b194656fc7aed7779ca3025a17/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

# --- de-de hyperlinks ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/68472668c37039b194656fc7aed7779ca3025a17/e2e/f29e047e-44c6-4de6-bfaa-fb03f56fc80b.md", [Type]::Missing, [Type]::Missing, "ffffcafed8d6-ffc7-4f2e-b2a7-3fd220cb73fa.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/40889a4ebe8bfd097568b4cc30e75c35cbfa571f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f29e047e-44c6-4de6-bfaa-fb03f56fc80b.543a0a31b72e48fbef1d9e1c94867f8b9fb36ec5.de-de.xlf", [Type]::Missing, [Type]::Missing, "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/4a62a6571c575aab2adb205669a75fb379d6b398/e2e/f29e047e-44c6-4de6-bfaa-fb03f56fc80b.md", [Type]::Missing, [Type]::Missing, "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/56392af51a58175f112d48fd6ea2b171fc9f939c/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f29e047e-44c6-4de6-bfaa-fb03f56fc80b.543a0a31b72e48fbef1d9e1c94867f8b9fb36ec5.de-de.xlf", [Type]::Missing, [Type]::Missing, "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/68472668c37039b194656fc7aed7779ca3025a17/e2e/ffffcafed8d6-ffc7-4f2e-b2a7-3fd220cb73fa.md", [Type]::Missing, [Type]::Missing, "ffffff5d935b58-0829-45f7-9e5b-8e53a6fa3f9f.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a8a93091275577b9bfbd16594fef9af79cc9368d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.de-de.xlf", [Type]::Missing, [Type]::Missing, "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/47b5b90b9b93c8dabd8651565b0a464410664acf/e2e/a28591f3-d3b4-45d5-86e2-cebe1e59fd36.md", [Type]::Missing, [Type]::Missing, "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/90cc793eeed73eb390e6b75487420be6d0e5f10b/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.de-de.xlf", [Type]::Missing, [Type]::Missing, "a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/68472668c37039b194656fc7aed7779ca3025a17/e2e/ffffff5d935b58-0829-45f7-9e5b-8e53a6fa3f9f.md", [Type]::Missing, [Type]::Missing, "f29e047e-44c6-4de6-bfaa-fb03f56fc80b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a8a93091275577b9bfbd16594fef9af79cc9368d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.de-de.xlf", [Type]::Missing, [Type]::Missing, "f29e047e-44c6-4de6-bfaa-fb03f56fc80b.543a0a31b72e48fbef1d9e1c94867f8b9fb36ec5.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/47b5b90b9b93c8dabd8651565b0a464410664acf/e2e/a28591f3-d3b4-45d5-86e2-cebe1e59fd36.md", [Type]::Missing, [Type]::Missing, "f29e047e-44c6-4de6-bfaa-fb03f56fc80b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/90cc793eeed73eb390e6b75487420be6d0e5f10b/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a28591f3-d3b4-45d5-86e2-cebe1e59fd36.c3fbac93e521f111b4aaa57a54c4db2453a3aeef.de-de.xlf", [Type]::Missing, [Type]::Missing, "f29e047e-44c6-4de6-bfaa-fb03f56fc80b.543a0a31b72e48fbef1d9e1c94867f8b9fb36ec5.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/68472668c37039b194656fc7aed7779ca3025a17/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null
